$wb = $excel.ActiveWorkbook

# Update "展览" sheet (Exhibition) - F3 and F4 "想去人数" (want-to-go count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 113
$ws1.Range("F4").Value = 10

# Update "全部类型" sheet (All Types) - same rows mirrored here
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 113
$ws4.Range("F4").Value = 10
